$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Static_Path" header to include an example path, and the
# sample static path value in the row below it.
$ws.Range("E1").Value = "Static_Path 1/101/1/13"
$ws.Range("E2").Value = "1/201/1/8"

# Column E needs to grow a bit to fit the longer header/value
# (COM ColumnWidth differs from the stored OOXML width by 5/6 of a
# character for this sheet's font, so subtract that offset to land on
# the desired stored width of 23.5).
$ws.Columns.Item(5).ColumnWidth = 22.666666666666668

# Move the active selection from F15 to F2
$ws.Range("F2").Select()
